$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.366.03"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "1.846.32"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'240.64"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").Value = "'0.6276"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'0.07488"
$ws.Range("E8").Value = "  -2.75%  "

$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").Value = "'24.38"
$ws.Range("E10").Value = "  -2.51%  "

$ws.Range("D11").Value = "'0.07738"
$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("D12").Value = "1.846.56"
$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("D13").Value = "'5.003"
$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("D14").Value = "'0.6789"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "'0.00001030"
$ws.Range("E15").Value = "  -5.25%  "

$ws.Range("D16").Value = "'82.64"
$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "2.110.51"
$ws.Range("E17").Value = "  -0.94%  "

$ws.Range("D18").Value = "'6.093"
$ws.Range("E18").Value = "  -1.58%  "

$ws.Range("D19").Value = "29.398.24"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").Value = "'12.29"
$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "'7.422"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").Value = "'159.00"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("D26").Value = "'0.1385"
$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("D27").Value = "'8.392"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D28").Value = "'17.57"
$ws.Range("E28").Value = "  -0.80%  "

$ws.Range("D29").Value = "'1.398"
$ws.Range("E29").Value = "  +3.29%  "

$ws.Range("D30").Value = "'1.476"
$ws.Range("E30").Value = "  +1.10%  "

$ws.Range("D31").Value = "'0.05701"
$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("D32").Value = "'4.110"
$ws.Range("E32").Value = "  -0.32%  "

$ws.Range("D33").Value = "'4.050"

$ws.Range("D34").Value = "'1.819"
$ws.Range("E34").Value = "  -1.38%  "

$ws.Range("D35").Value = "'1.149"
$ws.Range("E35").Value = "  -1.31%  "

$ws.Range("D36").Value = "'0.6930"
$ws.Range("E36").Value = "  -1.86%  "

$ws.Range("D37").Value = "'2.587"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").Value = "'2.851"
$ws.Range("E38").Value = "  +3.58%  "

$ws.Range("D39").Value = "1.250.91"
$ws.Range("E39").Value = "  +2.09%  "

$ws.Range("D40").Value = "'0.01816"
$ws.Range("E40").Value = "  +1.49%  "

$ws.Range("D41").Value = "'6.506"
$ws.Range("E41").Value = "  +0.88%  "

$ws.Range("D42").Value = "'0.9058"
$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("E44").Value = "  -1.21%  "

$ws.Range("D45").Value = "'101.30"
$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("D46").Value = "'65.80"
$ws.Range("E46").Value = "  -0.41%  "

$ws.Range("D47").Value = "'7.079"
$ws.Range("E47").Value = "  -1.46%  "

$ws.Range("D48").Value = "'0.1161"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.960"
$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.3938"
$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.00000000114"
$ws.Range("E51").Value = "  -4.66%  "
